# 10-May-2017 MainPO, Reclaims, CD Fee will be modified by now
$wb = $excel.ActiveWorkbook

# "acceptQueue": unit number correction (336985 -> 366965), keep as quoted text
$wsAcceptQueue = $wb.Worksheets.Item("acceptQueue")
$wsAcceptQueue.Range("A2").Value = "'366965"

# "Unit to Reconcile": clear the stale reconciled-unit value, keep formatting
$wsUnitToReconcile = $wb.Worksheets.Item("Unit to Reconcile")
$wsUnitToReconcile.Range("A2").ClearContents()

# "PO_Detail": unit number correction (00997989 -> 00998376), drop the quote-prefix style
$wsPoDetail = $wb.Worksheets.Item("PO_Detail")
$wsPoDetail.Range("A2").Value = "'00998376"
$wsPoDetail.Range("A2").Style = "Normal"

# "Unit_to_Reconcile_Output": clear out last run's output row
$wsOutput = $wb.Worksheets.Item("Unit_to_Reconcile_Output")
$wsOutput.Rows(2).ClearContents()

# Restore per-sheet selections (select acceptQueue's cell first so the final
# active sheet/tab stays on PO_Detail, matching the original activeTab)
$wsAcceptQueue.Range("A3").Select()
$wsPoDetail.Range("A2").Select()
